$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "748×7=" "229×3="
Replace-Text "225×8=" "401×4="
Replace-Text "798×9=" "490×8="
Replace-Text "391×4=" "419×6="
Replace-Text "854×4=" "716×4="
Replace-Text "145×3=" "178×9="
Replace-Text "292×9=" "345×7="
Replace-Text "219×9=" "486×8="
Replace-Text "346×6=" "226×3="
Replace-Text "381×3=" "843×3="
Replace-Text "113×6=" "971×8="
Replace-Text "392×4=" "931×9="
Replace-Text "523×6=" "913×3="
Replace-Text "690×4=" "925×2="
Replace-Text "951×7=" "284×3="
Replace-Text "234×8=" "623×8="
Replace-Text "236×8=" "803×9="
Replace-Text "140×7=" "920×7="
Replace-Text "746×2=" "540×3="
Replace-Text "894×6=" "820×2="
Replace-Text "594×5=" "386×4="
Replace-Text "399×8=" "878×5="
Replace-Text "661×3=" "514×8="
Replace-Text "168×3=" "424×7="
Replace-Text "364×2=" "216×9="

Write-Host "Done"
